$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp text in A1
$ws.Range("A1").Value = 'Datos actualizados a 21 de Marzo de 2020 a las 17:16'

# Rows 7-9
$arr = New-Object 'object[,]' 3,8
$arr[0,0] = 'Estados Unidos'
$arr[0,1] = 22019
$arr[0,2] = 2636
$arr[0,3] = 147
$arr[0,4] = 21591
$arr[0,5] = 64
$arr[0,6] = 25
$arr[0,7] = 281
$arr[1,0] = 'Alemania'
$arr[1,1] = 21682
$arr[1,2] = 1834
$arr[1,3] = 209
$arr[1,4] = 21398
$arr[1,5] = 2
$arr[1,6] = 7
$arr[1,7] = 75
$arr[2,0] = 'Iran'
$arr[2,1] = 20610
$arr[2,2] = 966
$arr[2,3] = 7635
$arr[2,4] = 11419
$arr[2,5] = 0
$arr[2,6] = 123
$arr[2,7] = 1556
$ws.Range("A7:H9").Value = $arr

# Rows 14-14
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Paises Bajos'
$arr[0,1] = 3631
$arr[0,2] = 637
$arr[0,3] = 2
$arr[0,4] = 3493
$arr[0,5] = 354
$arr[0,6] = 30
$arr[0,7] = 136
$ws.Range("A14:H14").Value = $arr

# Rows 17-17
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Noruega'
$arr[0,1] = 2050
$arr[0,2] = 91
$arr[0,3] = 1
$arr[0,4] = 2042
$arr[0,5] = 28
$arr[0,6] = 0
$arr[0,7] = 7
$ws.Range("A17:H17").Value = $arr

# Rows 34-36
$arr = New-Object 'object[,]' 3,8
$arr[0,0] = 'Grecia'
$arr[0,1] = 530
$arr[0,2] = 35
$arr[0,3] = 19
$arr[0,4] = 498
$arr[0,5] = 20
$arr[0,6] = 3
$arr[0,7] = 13
$arr[1,0] = 'Finlandia'
$arr[1,1] = 521
$arr[1,2] = 71
$arr[1,3] = 10
$arr[1,4] = 510
$arr[1,5] = 2
$arr[1,6] = 1
$arr[1,7] = 1
$arr[2,0] = 'Islandia'
$arr[2,1] = 473
$arr[2,2] = 64
$arr[2,3] = 5
$arr[2,4] = 467
$arr[2,5] = 1
$arr[2,6] = 1
$arr[2,7] = 1
$ws.Range("A34:H36").Value = $arr

# Rows 46-51
$arr = New-Object 'object[,]' 6,8
$arr[0,0] = 'India'
$arr[0,1] = 321
$arr[0,2] = 72
$arr[0,3] = 23
$arr[0,4] = 293
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 5
$arr[1,0] = 'Filipinas'
$arr[1,1] = 307
$arr[1,2] = 77
$arr[1,3] = 13
$arr[1,4] = 275
$arr[1,5] = 1
$arr[1,6] = 1
$arr[1,7] = 19
$arr[2,0] = 'Estonia'
$arr[2,1] = 306
$arr[2,2] = 23
$arr[2,3] = 2
$arr[2,4] = 304
$arr[2,5] = 1
$arr[2,6] = 0
$arr[2,7] = 0
$arr[3,0] = 'Rusia'
$arr[3,1] = 306
$arr[3,2] = 53
$arr[3,3] = 16
$arr[3,4] = 289
$arr[3,5] = 0
$arr[3,6] = 0
$arr[3,7] = 1
$arr[4,0] = 'Barein'
$arr[4,1] = 305
$arr[4,2] = 7
$arr[4,3] = 125
$arr[4,4] = 179
$arr[4,5] = 4
$arr[4,6] = 0
$arr[4,7] = 1
$arr[5,0] = 'Egipto'
$arr[5,1] = 285
$arr[5,2] = 0
$arr[5,3] = 42
$arr[5,4] = 235
$arr[5,5] = 0
$arr[5,6] = 0
$arr[5,7] = 8
$ws.Range("A46:H51").Value = $arr

# Rows 60-63
$arr = New-Object 'object[,]' 4,8
$arr[0,0] = 'Colombia'
$arr[0,1] = 196
$arr[0,2] = 51
$arr[0,3] = 1
$arr[0,4] = 195
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[1,0] = 'Kuwait'
$arr[1,1] = 176
$arr[1,2] = 17
$arr[1,3] = 27
$arr[1,4] = 149
$arr[1,5] = 5
$arr[1,6] = 0
$arr[1,7] = 0
$arr[2,0] = 'Bulgaria'
$arr[2,1] = 163
$arr[2,2] = 36
$arr[2,3] = 3
$arr[2,4] = 157
$arr[2,5] = 3
$arr[2,6] = 0
$arr[2,7] = 3
$arr[3,0] = 'Armenia'
$arr[3,1] = 160
$arr[3,2] = 24
$arr[3,3] = 1
$arr[3,4] = 159
$arr[3,5] = 2
$arr[3,6] = 0
$arr[3,7] = 0
$ws.Range("A60:H63").Value = $arr

# Rows 75-75
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Argelia'
$arr[0,1] = 102
$arr[0,2] = 8
$arr[0,3] = 43
$arr[0,4] = 47
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = 12
$ws.Range("A75:H75").Value = $arr

# Rows 77-77
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Bosnia y Herzegovina'
$arr[0,1] = 93
$arr[0,2] = 4
$arr[0,3] = 2
$arr[0,4] = 90
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 1
$ws.Range("A77:H77").Value = $arr

# Rows 83-89
$arr = New-Object 'object[,]' 7,8
$arr[0,0] = 'Republica de Chipre'
$arr[0,1] = 84
$arr[0,2] = 9
$arr[0,3] = 0
$arr[0,4] = 84
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 0
$arr[1,0] = 'Brunei'
$arr[1,1] = 83
$arr[1,2] = 5
$arr[1,3] = 1
$arr[1,4] = 82
$arr[1,5] = 2
$arr[1,6] = 0
$arr[1,7] = 0
$arr[2,0] = 'Lituania'
$arr[2,1] = 83
$arr[2,2] = 14
$arr[2,3] = 1
$arr[2,4] = 81
$arr[2,5] = 1
$arr[2,6] = 0
$arr[2,7] = 1
$arr[3,0] = 'Moldavia'
$arr[3,1] = 80
$arr[3,2] = 14
$arr[3,3] = 1
$arr[3,4] = 78
$arr[3,5] = 3
$arr[3,6] = 0
$arr[3,7] = 1
$arr[4,0] = 'Sri Lanka'
$arr[4,1] = 77
$arr[4,2] = 4
$arr[4,3] = 3
$arr[4,4] = 74
$arr[4,5] = 2
$arr[4,6] = 0
$arr[4,7] = 0
$arr[5,0] = 'Albania'
$arr[5,1] = 76
$arr[5,2] = 6
$arr[5,3] = 2
$arr[5,4] = 72
$arr[5,5] = 2
$arr[5,6] = 0
$arr[5,7] = 2
$arr[6,0] = 'Bielorrusia'
$arr[6,1] = 76
$arr[6,2] = 7
$arr[6,3] = 15
$arr[6,4] = 61
$arr[6,5] = 0
$arr[6,6] = 0
$arr[6,7] = 0
$ws.Range("A83:H89").Value = $arr

# Rows 95-95
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Camboya'
$arr[0,1] = 53
$arr[0,2] = 2
$arr[0,3] = 2
$arr[0,4] = 51
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A95:H95").Value = $arr
